$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 30
$ws.Range("B3").Value = 3.631442785263062
$ws.Range("B4").Value = 0.5939148664474487
$ws.Range("B5").Value = 2.359014749526978
$ws.Range("B6").Value = 3.315656900405884
$ws.Range("B7").Value = 3.535092830657959
$ws.Range("B8").Value = 3.891824781894684
$ws.Range("B9").Value = 5.353085517883301
